$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append six new (French-localized) master-data rows after the existing 13 rows.
# id, name, descr, lang_code ; is_active/cr_by/cr_dtimes follow the same pattern
# as every other row in the sheet (TRUE / superadmin / now()).
$newRows = @(
    @(10013, "Pré-inscription",        "Portail Web pour les pré-inscriptions"),
    @(10014, "Client dinscription",    "Application de bureau pour les inscriptions"),
    @(10015, "Processeur dinscription","Demande de post-inscription"),
    @(10016, "Authentification ID",    "Application pour lauthentification du fournisseur de services tiers"),
    @(10017, "Contrôle didentité",     "Portail Web pour la configuration dapplications"),
    @(10018, "Portail Résident",       "Portail Web pour les services de génération de post-ID")
)

$r = 14
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = "fra"
    $ws.Cells.Item($r, 5).Value = $true
    $ws.Cells.Item($r, 6).Value = "superadmin"
    $ws.Cells.Item($r, 7).Value = "now()"
    $r++
}

# Column width tweaks (id column narrower/best-fit, name column a bit wider).
$ws.Columns.Item(1).ColumnWidth = 5
$ws.Columns.Item(2).ColumnWidth = 18.5

# Scroll the view down and select the remainder of the sheet below the data
# (mirrors the author having clicked the row-20 header after typing the data).
$ws.Range("A20:XFD1048576").Select() | Out-Null

# Page setup: Letter/A4-class "9" paper (A4), portrait orientation.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
